$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 1344.2222
$ws.Range("J17").Value = 1344.2222
$ws.Range("L17").Value = 4032.6666
$ws.Range("N17").Value = -4368.6666

# Row 19
$ws.Range("H19").Value = 2333
$ws.Range("I19").Value = 1499.75
$ws.Range("K19").Value = 1499.75
$ws.Range("M19").Value = -1324.75

# Row 32
$ws.Range("H32").Value = 0
$ws.Range("J32").Value = 0
$ws.Range("L32").Value = 0
$ws.Range("N32").ClearContents()

# Row 40
$ws.Range("H40").Value = 3914.3809
$ws.Range("I40").Value = 2677.5557
$ws.Range("K40").Value = 2677.5557
$ws.Range("M40").Value = -2502.5557

# Row 51
$ws.Range("H51").Value = 16155.857
$ws.Range("J51").Value = 5918.4
$ws.Range("L51").Value = 5918.4
$ws.Range("N51").Value = -6886.4

# Row 86
$ws.Range("H86").Value = 47842856
$ws.Range("I86").Value = 56189436
$ws.Range("K86").Value = 56189436
$ws.Range("M86").Value = -56188313

# Row 89
$ws.Range("H89").Value = 47842856
$ws.Range("I89").Value = 56189436
$ws.Range("K89").Value = 280947180
$ws.Range("M89").Value = -280941564

# Row 116
$ws.Range("H116").Value = 41674830
$ws.Range("J116").Value = 14333.333
$ws.Range("L116").Value = 14333.333
$ws.Range("N116").Value = -21217.333

# Row 132
$ws.Range("H132").Value = 1203.8959
$ws.Range("I132").Value = 1186.9574
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 3560.8722
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -1030.8722
$ws.Range("N132").Value = -11060

# Row 137
$ws.Range("H137").Value = 2799
$ws.Range("I137").Value = 4577.4
$ws.Range("J137").Value = 2025.7826
$ws.Range("K137").Value = 13732.2
$ws.Range("L137").Value = 6077.3478
$ws.Range("M137").Value = -11182.2
$ws.Range("N137").Value = -11177.3478

# Row 138
$ws.Range("H138").Value = 6044.3696
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 6044.3696
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 18133.1088
$ws.Range("M138").ClearContents()
$ws.Range("N138").Value = -28413.1088

$ws = $wb.Worksheets.Item("ARM")
# Row 57
$ws.Range("H57").Value = 4864
$ws.Range("I57").Value = 4864
$ws.Range("K57").Value = 4864
$ws.Range("M57").Value = -4380

# Row 61
$ws.Range("H61").Value = 10469.385
$ws.Range("I61").Value = 2549.25
$ws.Range("J61").Value = 13989.444
$ws.Range("K61").Value = 2549.25
$ws.Range("L61").Value = 13989.444
$ws.Range("M61").Value = -2337.25
$ws.Range("N61").Value = -14413.444

# Row 122
$ws.Range("H122").Value = 6983.364
$ws.Range("I122").Value = 4184.6
$ws.Range("K122").Value = 12553.8
$ws.Range("M122").Value = -10103.8

# Row 126
$ws.Range("H126").Value = 5267.5
$ws.Range("I126").Value = 5267.5
$ws.Range("K126").Value = 15802.5
$ws.Range("M126").Value = -13332.5

# Row 132
$ws.Range("H132").Value = 7206.641
$ws.Range("I132").Value = 6242.0527
$ws.Range("K132").Value = 18726.1581
$ws.Range("M132").Value = -16196.1581

# Row 136
$ws.Range("H136").Value = 10469.385
$ws.Range("I136").Value = 2549.25
$ws.Range("J136").Value = 13989.444
$ws.Range("K136").Value = 7647.75
$ws.Range("L136").Value = 41968.33199999999
$ws.Range("M136").Value = -5097.75
$ws.Range("N136").Value = -47068.33199999999

$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 955263.25
$ws.Range("I105").Value = 1430466.2
$ws.Range("K105").Value = 1430466.2
$ws.Range("M105").Value = -1428719.2

# Row 113
$ws.Range("H113").Value = 5209.4546
$ws.Range("I113").Value = 5209.4546
$ws.Range("K113").Value = 5209.4546
$ws.Range("M113").Value = -3039.4546

# Row 134
$ws.Range("H134").Value = 6445.452
$ws.Range("I134").Value = 2606.0667
$ws.Range("J134").Value = 8578.444
$ws.Range("K134").Value = 7818.2001
$ws.Range("L134").Value = 25735.332
$ws.Range("M134").Value = -5283.2001
$ws.Range("N134").Value = -30805.332

$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 53.909092
$ws.Range("I7").Value = 62.5
$ws.Range("K7").Value = 62.5
$ws.Range("M7").Value = 50.5

# Row 31
$ws.Range("H31").Value = 5850.6855
$ws.Range("I31").Value = 2525.1135
$ws.Range("J31").Value = 11478.577
$ws.Range("K31").Value = 2525.1135
$ws.Range("L31").Value = 11478.577
$ws.Range("M31").Value = -2230.1135
$ws.Range("N31").Value = -12068.577

# Row 34
$ws.Range("H34").Value = 5850.6855
$ws.Range("I34").Value = 2525.1135
$ws.Range("J34").Value = 11478.577
$ws.Range("K34").Value = 2525.1135
$ws.Range("L34").Value = 11478.577
$ws.Range("M34").Value = -2323.1135
$ws.Range("N34").Value = -11882.577

# Row 107
$ws.Range("H107").Value = 1947.8572
$ws.Range("I107").Value = 427.7
$ws.Range("K107").Value = 427.7
$ws.Range("M107").Value = 1492.3

# Row 132
$ws.Range("H132").Value = 3649.2104
$ws.Range("I132").Value = 1753.0892
$ws.Range("K132").Value = 5259.267599999999
$ws.Range("M132").Value = -2729.267599999999

# Row 134
$ws.Range("H134").Value = 4374.127
$ws.Range("I134").Value = 1933.762
$ws.Range("J134").Value = 9254.857
$ws.Range("K134").Value = 5801.286
$ws.Range("L134").Value = 27764.571
$ws.Range("M134").Value = -3266.286
$ws.Range("N134").Value = -32834.571

$ws = $wb.Worksheets.Item("CUL")
# Row 23
$ws.Range("H23").Value = 217.96
$ws.Range("I23").Value = 170.9
$ws.Range("J23").Value = 249.33333
$ws.Range("K23").Value = 512.7
$ws.Range("L23").Value = 747.99999
$ws.Range("M23").Value = -277.7
$ws.Range("N23").Value = -1217.99999

# Row 121
$ws.Range("H121").Value = 3847327
$ws.Range("J121").Value = 4167912.5
$ws.Range("L121").Value = 12503737.5
$ws.Range("N121").Value = -12506357.5

$ws = $wb.Worksheets.Item("GSM")
# Row 80
$ws.Range("H80").Value = 2945.8
$ws.Range("I80").Value = 3058.25
$ws.Range("J80").Value = 2496
$ws.Range("K80").Value = 3058.25
$ws.Range("L80").Value = 2496
$ws.Range("M80").Value = -2060.25
$ws.Range("N80").Value = -4492

# Row 83
$ws.Range("H83").Value = 2945.8
$ws.Range("I83").Value = 3058.25
$ws.Range("J83").Value = 2496
$ws.Range("K83").Value = 15291.25
$ws.Range("L83").Value = 12480
$ws.Range("M83").Value = -10299.25
$ws.Range("N83").Value = -22464

# Row 113
$ws.Range("H113").Value = 250264.05
$ws.Range("I113").Value = 669536.75
$ws.Range("K113").Value = 669536.75
$ws.Range("M113").Value = -667366.75

# Row 126
$ws.Range("H126").Value = 3519.1538
$ws.Range("J126").Value = 4328.2856
$ws.Range("L126").Value = 12984.8568
$ws.Range("N126").Value = -17924.8568

# Row 132
$ws.Range("H132").Value = 4662.51
$ws.Range("I132").Value = 2289.5
$ws.Range("K132").Value = 6868.5
$ws.Range("M132").Value = -4338.5

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 6499.3687
$ws.Range("I7").Value = 5004
$ws.Range("K7").Value = 5004
$ws.Range("M7").Value = -4892

# Row 16
$ws.Range("H16").Value = 819.5714
$ws.Range("I16").Value = 872.8333
$ws.Range("K16").Value = 872.8333
$ws.Range("M16").Value = -702.8333

# Row 46
$ws.Range("H46").Value = 7940609.5
$ws.Range("J46").Value = 8551310
$ws.Range("L46").Value = 8551310
$ws.Range("N46").Value = -8551686

# Row 61
$ws.Range("H61").Value = 6457.6665
$ws.Range("I61").Value = 1530
$ws.Range("J61").Value = 7689.5835
$ws.Range("K61").Value = 1530
$ws.Range("L61").Value = 7689.5835
$ws.Range("M61").Value = -1328
$ws.Range("N61").Value = -8093.5835

# Row 82
$ws.Range("H82").Value = 3046.1667
$ws.Range("I82").Value = 3020.6
$ws.Range("K82").Value = 3020.6
$ws.Range("M82").Value = -2659.6

# Row 85
$ws.Range("H85").Value = 3046.1667
$ws.Range("I85").Value = 3020.6
$ws.Range("K85").Value = 3020.6
$ws.Range("M85").Value = -1772.6

# Row 113
$ws.Range("H113").Value = 6457.6665
$ws.Range("I113").Value = 1530
$ws.Range("J113").Value = 7689.5835
$ws.Range("K113").Value = 1530
$ws.Range("L113").Value = 7689.5835
$ws.Range("M113").Value = 640
$ws.Range("N113").Value = -12029.5835

# Row 126
$ws.Range("H126").Value = 6499.3687
$ws.Range("I126").Value = 5004
$ws.Range("K126").Value = 15012
$ws.Range("M126").Value = -12542

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 227481.89
$ws.Range("I122").Value = 503288.38
$ws.Range("K122").Value = 1509865.14
$ws.Range("M122").Value = -1507415.14

# Row 126
$ws.Range("H126").Value = 3201.2307
$ws.Range("I126").Value = 2311.6
$ws.Range("K126").Value = 6934.799999999999
$ws.Range("M126").Value = -4464.799999999999

# Row 136
$ws.Range("H136").Value = 28604408
$ws.Range("I136").Value = 55556636
$ws.Range("K136").Value = 166669908
$ws.Range("M136").Value = -166667358

